# Monte Carlo Type AI - update the two HP-tracking columns (H = HP CUR,
# I = HP MAX) for both Pokemon rows to reflect the new simulated values,
# and move the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Pokemon #321): HP CUR / HP MAX 200 -> 69
$ws.Range("H2").Value = 69
$ws.Range("I2").Value = 69

# Row 3 (Pokemon #1237): HP CUR / HP MAX 189 -> 102
$ws.Range("H3").Value = 102
$ws.Range("I3").Value = 102

# Move the active selection from Z3 to H3
$ws.Activate()
$ws.Range("H3").Select()
